$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C19").Value = 7639
$ws.Range("C20:C59").Value = 7312
$ws.Range("C92:C166").Value = 7310
$ws.Range("C167:C171").Value = 7295
